$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to append: dates in column A, times in column B (rows 7-15)
$dates = @(45187, 45190, 45191, 45194, 45195, 45196, 45197, 45198, 45199)
$times = @(0.40625, 0.40902777777777777, 0.42708333333333331, 0.4201388888888889, 0.4236111111111111, 0.4236111111111111, 0.4201388888888889, 0.4236111111111111, 0.43055555555555558)

$row = 7
for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $times[$i]
    $row++
}

# Reuse the existing formatting from rows 2-6 (date / time number formats +
# centered, wrapped alignment) by copying it onto the freshly written rows.
$ws.Range("A2:B6").Copy()
$ws.Range("A7:B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection like in the diff
$ws.Range("G17").Select()
